# Update the Pcf reporting template:
#  - widen column A on the "semantic_aspect_model_schema" sheet
#  - rename the "id" field display name to "dtwin_id" (schema + description sheets)
#  - refresh the generation metadata (commit hash / url / date)

$wb = $excel.ActiveWorkbook

$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")
$wsDescription = $wb.Worksheets.Item("description")
$wsMetadata = $wb.Worksheets.Item("metadata")

# Widen column A (was 2.4 characters wide, target 9.6) so the "dtwin_id"
# header name doesn't collide with the aspect-model's own "id" column names.
# (ColumnWidth is snapped to the host's pixel grid, so we feed it the value
# that resolves to the closest on-grid width to 9.6.)
$wsSchema.Columns.Item(1).ColumnWidth = (53.0/6.0)

# Rename the digital-twin id field's display name on both sheets.
$wsSchema.Range("A1").Value = "dtwin_id"
$wsDescription.Range("A5").Value = "dtwin_id"

# Record the commit this reporting template was regenerated from.
$wsMetadata.Range("B2").Value = "41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B3").Value = "https://github.com/dataspacesolutions/sldt-semantic-models/commit/41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B4").Value = "2025-03-10 14:48:29+00:00"
